$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add human-readable descriptions of computer-generated masks compared to ground truth
$ws.Range('J40').Value = 'decent - detecting vessels'
$ws.Range('K41').Value = 'TERRIBLE - juxtapleural'
$ws.Range('K41').Style = 'Bad'
$ws.Range('Q41').Value = 'ideas: if elongation is close to 1, delete object'
$ws.Range('K42').Value = 'boundary is bad - detecting lung edge'
$ws.Range('K43').Value = 'boundary is bad - detecting lung edge'
$ws.Range('K44').Value = 'terrible - juxtapleural'
$ws.Range('K44').Style = 'Bad'
$ws.Range('K45').Value = 'terrible - juxtapleural'
$ws.Range('K45').Style = 'Bad'
$ws.Range('M45').Value = 'use in PPT - almost juxtapleural'
$ws.Range('K46').Value = 'terrible - juxtapleural'
$ws.Range('K46').Style = 'Bad'
$ws.Range('M46').Value = 'use in PPT - almost juxtapleural'
$ws.Range('K47').Value = 'terrible - juxtapleural'
$ws.Range('K47').Style = 'Bad'
$ws.Range('J48').Value = 'decent - detecting vessels'
$ws.Range('J49').Value = 'decent - detecting vessels'
$ws.Range('K50').Value = 'bad - detecting vessels and lung edge'
$ws.Range('K50').Style = 'Bad'
$ws.Range('K51').Value = 'boundary is bad - detecting lung edge'
$ws.Range('J52').Value = 'decent  '
$ws.Range('J53').Value = 'decent  small vessel detected'
$ws.Range('J54').Value = 'decent - lung edge detected despite being juxtapleural'
$ws.Range('J55').Value = 'decent - lung edge detected despite being juxtapleural'
$ws.Range('J56').Value = 'decent - lung edge detected despite being juxtapleural'
$ws.Range('K56').Value = 'floodfill and interior messed up'
$ws.Range('J57').Value = 'decent'
$ws.Range('J58').Value = 'decent - boundary detects lung edges'
$ws.Range('J59').Value = 'decent - boundary detects lung edges'
$ws.Range('K60').Value = 'TERRIBLE - juxtapleural'
$ws.Range('K60').Style = 'Bad'
$ws.Range('J61').Value = 'decent - detecting vessels'
$ws.Range('J62').Value = 'decent - detecting some vessels'
$ws.Range('J63').Value = 'decent - detecting some vessels'
$ws.Range('J64').Value = 'decent - detecting some vessels'
$ws.Range('K65').Value = 'TERRIBLE - juxtapleural'
$ws.Range('K65').Style = 'Bad'
$ws.Range('J66').Value = 'decent - detecting some vessels'
$ws.Range('J67').Value = 'decent - detecting some vessels'
$ws.Range('K68').Value = 'TERRIBLE - juxtapleural'
$ws.Range('K68').Style = 'Bad'
$ws.Range('M68').Value = 'fix the floodfill mask'
$ws.Range('J69').Value = 'decent - detecting lung edge'
$ws.Range('K70').Value = 'TERRIBLE - juxtapleural'
$ws.Range('K70').Style = 'Bad'
$ws.Range('K71').Value = 'TERRIBLE - confused with vessels'
$ws.Range('K71').Style = 'Bad'

# Update selection to reflect final active cell
$ws.Range('P63').Select()
